# Auto-generated edit script: applies cell-level changes per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44195
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 1500
$ws.Range("D3").Value = 45126
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 24000
$ws.Range("S3").Value = 2400
$ws.Range("D4").Value = 45126
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 2000
$ws.Range("D5").Value = 45126
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 1500
$ws.Range("D6").Value = 45126
$ws.Range("L6").Value = "Tercera"
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("S6").Value = 1200
$ws.Range("T6").Value = 10
$ws.Range("D7").Value = 44400
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 1500
$ws.Range("O7").Value = 1500
$ws.Range("P7").Value = 1500
$ws.Range("S7").Value = 1500
$ws.Range("D8").Value = 44292
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("S8").Value = 1400
$ws.Range("T8").Value = 10
$ws.Range("D9").Value = 44904
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("S9").Value = 1500
$ws.Range("T9").Value = 10
$ws.Range("D10").Value = 44904
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 1000
$ws.Range("D11").Value = 45125
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 55
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 24000
$ws.Range("S11").Value = 2400
$ws.Range("D12").Value = 45125
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("S12").Value = 2000
$ws.Range("D13").Value = 45125
$ws.Range("L13").Value = "Segunda"
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("S13").Value = 1500
$ws.Range("D14").Value = 45125
$ws.Range("L14").Value = "Tercera"
$ws.Range("M14").Value = 45
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("S14").Value = 1200
$ws.Range("D15").Value = 44371
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 1800
$ws.Range("O15").Value = 1800
$ws.Range("P15").Value = 1800
$ws.Range("Q15").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S15").Value = 1800
$ws.Range("T15").Value = 1
$ws.Range("D16").Value = 44371
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 1200
$ws.Range("O16").Value = 1200
$ws.Range("P16").Value = 1200
$ws.Range("Q16").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S16").Value = 1200
$ws.Range("T16").Value = 1
$ws.Range("D17").Value = 45118
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 140
$ws.Range("N17").Value = 24000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 24000
$ws.Range("S17").Value = 2400
$ws.Range("D18").Value = 45118
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("S18").Value = 2000
$ws.Range("D19").Value = 45118
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 180
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 15000
$ws.Range("P19").Value = 15000
$ws.Range("S19").Value = 1500
$ws.Range("D20").Value = 45118
$ws.Range("L20").Value = "Tercera"
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 12000
$ws.Range("S20").Value = 1200
$ws.Range("D25").Value = 44391
$ws.Range("M25").Value = 15
$ws.Range("D26").Value = 44391
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 20
$ws.Range("N26").Value = 1000
$ws.Range("O26").Value = 1000
$ws.Range("P26").Value = 1000
$ws.Range("S26").Value = 1000
$ws.Range("D27").Value = 44309
$ws.Range("N27").Value = 1600
$ws.Range("O27").Value = 1600
$ws.Range("P27").Value = 1600
$ws.Range("S27").Value = 1600
$ws.Range("D28").Value = 44336
$ws.Range("M28").Value = 10
$ws.Range("N28").Value = 1500
$ws.Range("O28").Value = 1500
$ws.Range("P28").Value = 1500
$ws.Range("Q28").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("T28").Value = 1
$ws.Range("D29").Value = 44880
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 20000
$ws.Range("S29").Value = 2000
$ws.Range("D30").Value = 44880
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 180
$ws.Range("N30").Value = 15000
$ws.Range("O30").Value = 15000
$ws.Range("P30").Value = 15000
$ws.Range("S30").Value = 1500
